$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values for Jogos_da_Semana_FlashScore_2025-05-12

# Row 2
$ws.Range("J2").Value = 1.05
$ws.Range("L2").Value = 1.37
$ws.Range("T2").Value = 8
$ws.Range("AF2").Value = 12

# Row 3
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 3.6
$ws.Range("J3").Value = 1.03
$ws.Range("L3").Value = 1.22
$ws.Range("N3").Value = 1.79
$ws.Range("O3").Value = 1.99
$ws.Range("P3").Value = 1.33
$ws.Range("Q3").Value = 3.25
$ws.Range("T3").Value = 8.5
$ws.Range("AA3").Value = 6.5
$ws.Range("AD3").Value = 151

# Row 4
$ws.Range("G4").Value = 3.75
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 1.8
$ws.Range("J4").Value = 1.03
$ws.Range("K4").Value = 9.75
$ws.Range("L4").Value = 1.18
$ws.Range("M4").Value = 4.5
$ws.Range("N4").Value = 1.55
$ws.Range("O4").Value = 2.35
$ws.Range("P4").Value = 1.29
$ws.Range("Q4").Value = 3.35
$ws.Range("R4").Value = 1.53
$ws.Range("S4").Value = 2.35
$ws.Range("T4").Value = 14
$ws.Range("U4").Value = 26
$ws.Range("V4").Value = 13.5
$ws.Range("W4").Value = 60
$ws.Range("X4").Value = 32
$ws.Range("Y4").Value = 32
$ws.Range("Z4").Value = 9.75
$ws.Range("AA4").Value = 8.5
$ws.Range("AB4").Value = 13.5
$ws.Range("AC4").Value = 45
$ws.Range("AD4").Value = 300
$ws.Range("AE4").Value = 9.5
$ws.Range("AF4").Value = 11.25
$ws.Range("AG4").Value = 8.75
$ws.Range("AH4").Value = 17.5
$ws.Range("AI4").Value = 13.5
$ws.Range("AJ4").Value = 21

# Row 5
$ws.Range("G5").Value = 1.3
$ws.Range("H5").Value = 5.25
$ws.Range("J5").Value = 1.05
$ws.Range("K5").Value = 11
$ws.Range("N5").Value = 2.03
$ws.Range("O5").Value = 1.78
$ws.Range("V5").Value = 9.5
$ws.Range("X5").Value = 13
$ws.Range("Z5").Value = 9
$ws.Range("AE5").Value = 19
$ws.Range("AF5").Value = 51
$ws.Range("AG5").Value = 34
$ws.Range("AI5").Value = 101
$ws.Range("AJ5").Value = 101

# Row 7
$ws.Range("J7").Value = 1.17
$ws.Range("K7").Value = 5

# Row 16
$ws.Range("G16").Value = 3.1
$ws.Range("I16").Value = 2.3
$ws.Range("M16").Value = 2.65
$ws.Range("T16").Value = 8.5
$ws.Range("U16").Value = 15.5
$ws.Range("V16").Value = 11
$ws.Range("X16").Value = 29
$ws.Range("AA16").Value = 5.9
$ws.Range("AE16").Value = 6.8
$ws.Range("AF16").Value = 10.5
$ws.Range("AH16").Value = 24

# Row 21
$ws.Range("G21").Value = 1.87
$ws.Range("R21").Value = 1.8
$ws.Range("S21").Value = 1.8

# Row 22
$ws.Range("I22").Value = 3.1
$ws.Range("R22").Value = 1.77
$ws.Range("S22").Value = 1.92
$ws.Range("W22").Value = 23
$ws.Range("X22").Value = 21
$ws.Range("AD22").Value = 251
$ws.Range("AI22").Value = 26
$ws.Range("AJ22").Value = 34

# Row 23
$ws.Range("G23").Value = 2.35
$ws.Range("J23").Value = 1.06
$ws.Range("K23").Value = 10
$ws.Range("R23").Value = 1.69

# Row 24
$ws.Range("G24").Value = 1.82
$ws.Range("K24").Value = 10
$ws.Range("R24").Value = 1.87
$ws.Range("S24").Value = 1.77

# Row 25
$ws.Range("R25").Value = 1.5
$ws.Range("S25").Value = 2.37

# Row 26
$ws.Range("G26").Value = 2.15
$ws.Range("H26").Value = 3
$ws.Range("I26").Value = 3.2
$ws.Range("N26").Value = 2.1
$ws.Range("O26").Value = 1.57
$ws.Range("P26").Value = 1.42
$ws.Range("Q26").Value = 2.35
$ws.Range("T26").Value = 5.7
$ws.Range("U26").Value = 8.25
$ws.Range("V26").Value = 7.6
$ws.Range("W26").Value = 17
$ws.Range("X26").Value = 15.5
$ws.Range("Y26").Value = 25
$ws.Range("Z26").Value = 7.6
$ws.Range("AA26").Value = 5.2
$ws.Range("AB26").Value = 12.5
$ws.Range("AC26").Value = 60
$ws.Range("AD26").Value = 450
$ws.Range("AE26").Value = 7.2
$ws.Range("AF26").Value = 13
$ws.Range("AG26").Value = 9.5
$ws.Range("AH26").Value = 32
$ws.Range("AI26").Value = 24
$ws.Range("AJ26").Value = 32

# Row 27
$ws.Range("G27").Value = 5.7
$ws.Range("H27").Value = 3.9
$ws.Range("I27").Value = 1.47
$ws.Range("N27").Value = 1.83
$ws.Range("O27").Value = 1.78
$ws.Range("P27").Value = 1.37
$ws.Range("Q27").Value = 2.5
$ws.Range("T27").Value = 11.5
$ws.Range("U27").Value = 26
$ws.Range("V27").Value = 15.5
$ws.Range("W27").Value = 80
$ws.Range("X27").Value = 50
$ws.Range("Y27").Value = 50
$ws.Range("Z27").Value = 9.75
$ws.Range("AA27").Value = 6.7
$ws.Range("AB27").Value = 16
$ws.Range("AC27").Value = 75
$ws.Range("AD27").Value = 500
$ws.Range("AE27").Value = 5.3
$ws.Range("AF27").Value = 5.6
$ws.Range("AG27").Value = 7.1
$ws.Range("AH27").Value = 8.25
$ws.Range("AI27").Value = 10.5
$ws.Range("AJ27").Value = 24

# Row 28
$ws.Range("G28").Value = 3.75
$ws.Range("H28").Value = 3.05
$ws.Range("I28").Value = 1.93
$ws.Range("N28").Value = 2.1
$ws.Range("O28").Value = 1.57
$ws.Range("P28").Value = 1.42
$ws.Range("Q28").Value = 2.35
$ws.Range("T28").Value = 7.8
$ws.Range("U28").Value = 15.5
$ws.Range("V28").Value = 10.75
$ws.Range("W28").Value = 45
$ws.Range("X28").Value = 30
$ws.Range("Y28").Value = 37
$ws.Range("Z28").Value = 7.5
$ws.Range("AA28").Value = 5.3
$ws.Range("AB28").Value = 13
$ws.Range("AC28").Value = 65
$ws.Range("AD28").Value = 500
$ws.Range("AE28").Value = 5.3
$ws.Range("AF28").Value = 7.2
$ws.Range("AG28").Value = 7.2
$ws.Range("AH28").Value = 13.5
$ws.Range("AI28").Value = 14
$ws.Range("AJ28").Value = 25

# Row 29
$ws.Range("P29").Value = 1.22

# Row 30
$ws.Range("P30").Value = 1.19

# Row 31
$ws.Range("G31").Value = 1.9
$ws.Range("AJ31").Value = 41

# Row 39
$ws.Range("G39").Value = 2.45
$ws.Range("H39").Value = 2.9
$ws.Range("I39").Value = 3
$ws.Range("L39").Value = 1.57
$ws.Range("M39").Value = 2.25
$ws.Range("N39").Value = 2.88
$ws.Range("O39").Value = 1.4
$ws.Range("P39").Value = 1.58
$ws.Range("Q39").Value = 2.2
$ws.Range("T39").Value = 6
$ws.Range("U39").Value = 10
$ws.Range("V39").Value = 11
$ws.Range("W39").Value = 26
$ws.Range("X39").Value = 26
$ws.Range("AE39").Value = 6.5
$ws.Range("AF39").Value = 13
$ws.Range("AH39").Value = 34

# Row 40
$ws.Range("N40").Value = 1.84
$ws.Range("O40").Value = 1.84
$ws.Range("P40").Value = 1.33

# Row 42
$ws.Range("G42").Value = 1.1
$ws.Range("H42").Value = 10
$ws.Range("I42").Value = 12
$ws.Range("K42").Value = 34
$ws.Range("P42").Value = 1.11
$ws.Range("W42").Value = 8
$ws.Range("AD42").Value = 900
$ws.Range("AE42").Value = 41

# Row 43
$ws.Range("G43").Value = 5.5
$ws.Range("H43").Value = 5
$ws.Range("I43").Value = 1.39
$ws.Range("J43").Value = 1.01
$ws.Range("K43").Value = 17
$ws.Range("P43").Value = 1.25
$ws.Range("R43").Value = 1.8
$ws.Range("S43").Value = 1.91
$ws.Range("T43").Value = 19
$ws.Range("U43").Value = 34
$ws.Range("W43").Value = 67
$ws.Range("AA43").Value = 10
$ws.Range("AB43").Value = 19
$ws.Range("AE43").Value = 8.5
$ws.Range("AF43").Value = 8
$ws.Range("AH43").Value = 10
$ws.Range("AJ43").Value = 23

# Row 55
$ws.Range("G55").Value = 1.73
$ws.Range("R55").Value = 2.1
$ws.Range("S55").Value = 1.67
$ws.Range("X55").Value = 15

# Row 56
$ws.Range("N56").Value = 2.08
$ws.Range("O56").Value = 1.73

# Row 57
$ws.Range("G57").Value = 1.7
$ws.Range("H57").Value = 3.8
$ws.Range("I57").Value = 4.5
$ws.Range("T57").Value = 8
$ws.Range("W57").Value = 13
$ws.Range("AB57").Value = 15
$ws.Range("AD57").Value = 201
$ws.Range("AE57").Value = 15
$ws.Range("AG57").Value = 15
$ws.Range("AH57").Value = 51
$ws.Range("AJ57").Value = 41
